$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You work for a fitness center that maintains a database containing exhaustive information about potential leads, including age, gender, address, previous gym membership purchases, and attendance to prior free events. Your daily target is to contact 60 new leads within two hours to introduce a new gym membership. After one hour today, you called only 15 prospects. Most did not pick up the phone. Those you spoke to did not qualify for the offer.  What should you do next?",
        "ques_type": 2,
        "options": [
            "Identify 45 warm leads in the database, such as those who previously purchased a membership or often attend free events.",
            "Quickly check in with your manager for suggestions on improving your approach.",
            "Review and adjust the script to determine whether the lead is qualified for the offer, with a target of less than two minutes per lead.",
            "Analyze why the leads did not qualify for the offer to better select the next 45 prospects."
        ],
        "score": "Identify 45 warm leads in the database, such as those who previously purchased a membership or often attend free events."
    },
    {
        "title": "You work for a life coach who proposes high-ticket coaching programs to women on \u201chow to discover and live their life\u2019s purpose,\u201d with an entry price of $7,000 per program. Women interested in joining the program complete an application form and receive a follow-up call to assess whether they qualify. During one of your follow-up calls, a woman explains, \u201cI\u2019m very interested in joining, but I still need to discuss it with my husband.\u201d How should you respond?",
        "ques_type": 2,
        "options": [
            "Thank her for the time she spent with you and ask her to contact you once she has discussed it with her husband.",
            "Ask for a follow-up call with her and her husband to share the program\u2019s benefits directly with him.",
            "Ask if there is anything you can do to support her as she talks to her husband, and book another follow-up call with her.",
            "Highlight more benefits of the program, and book another follow-up call with her."
        ],
        "score": "Ask if there is anything you can do to support her as she talks to her husband, and book another follow-up call with her."
    },
    {
        "title": "You work for a company that provides transport services to individuals for fixed, non-negotiable prices and free additional services to increase customer loyalty. You receive a call from a client needing a last-minute driver and VIP car for a trip to a party on the same day. You give her the quote, and she says, \u201cIt\u2019s too expensive. Are you able to reduce the price?\u201d What should you say next to move the sale forward?",
        "ques_type": 2,
        "options": [
            "Suggest a standard car rather than a VIP one, and state the new price quote.",
            "Politely explain that the rate is non-negotiable, but emphasize that the company is highly rated, so the service is worth the rate. ",
            "Politely explain that the rate is non-negotiable, and offer her a free additional service.",
            "Ask more questions about her budget and what is most important for her during this trip."
        ],
        "score": "Ask more questions about her budget and what is most important for her during this trip."
    },
    {
        "title": "You work for a company selling online educational children\u2019s programs to parents. Today, you are calling an existing database of parents to introduce them to a new children\u2019s program the company has just launched. A man answers your first call. After introducing yourself, the client angrily interrupts you, explains that he is not happy with the current program he uses, and asks for a refund. You are not familiar with this program. How can you best react to this client?",
        "ques_type": 2,
        "options": [
            "Acknowledge his point of view, and question him further on areas he felt the company did not deliver what he was expecting for his child.",
            "Assess his level of anger and decide whether to transfer the call to someone familiar with the program for which he wants a refund.",
            "Empathize with him about his previous negative experience, take his information to organize his refund, then put your focus on presenting your new offer.",
            "Thank him for his honesty, take his information to organize his refund, and make a note to call back once the refund is processed."
        ],
        "score": "Acknowledge his point of view, and question him further on areas he felt the company did not deliver what he was expecting for his child."
    }
]
'@

# Trim the trailing newline that the here-string literal adds.
$newText = $newText.TrimEnd("`r", "`n")

# The sheet used to hold two rows: A1 was a placeholder number (bold,
# bordered, centered style) and A2 held the long serialized-questions
# string as a shared string. Collapse this back down to a single row:
# delete row 2 (shifting nothing else up, it's the last row) and put the
# (reformatted) long string directly into A1 instead, with no special
# formatting.
$ws.Range("A2").Delete()

$ws.Range("A1").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = $newText

# Writing a long multi-line string makes Excel auto-expand the row height;
# auto-fit it back down so no explicit row height is persisted.
$ws.Rows.Item(1).AutoFit()
